# Slide 1 ("Title 1" placeholder, shape index 1):
#   - retitle "Project Analysis" -> "The Battle of Neighborhoods"
#   - drop the bold/white title styling in favor of a regular-weight,
#     dark (#1F1F1F) OpenSans run, clearing any inherited text shadow
#   - the placeholder has <a:spAutoFit/>, so its box height recalculates
#     on its own once the text/font change lands - no manual resize needed
# Also remove the trailing "24Slides" attribution picture from the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)

# --- Retitle the first run of the title placeholder ---
$tf = $shp.TextFrame
$tr = $tf.TextRange
$titleRun = $tr.Characters(1, 17)              # "Project Analysis" = 17 chars
$titleRun.Text = "The Battle of Neighborhoods"

# --- Restyle that same run (now 27 chars long) ---
$titleRun = $tr.Characters(1, 27)
$titleRun.Font.Bold = $false
$titleRun.Font.Italic = $false
$titleRun.Font.Name = "OpenSans"
$titleRun.Font.Color.RGB = 0x1F1F1F

# --- Explicitly clear any inherited effect (writes an empty <a:effectLst/>) ---
$tf2 = $shp.TextFrame2
$tr2 = $tf2.TextRange
$titleRun2 = $tr2.Characters(1, 27)
$titleRun2.Font.Shadow = $false

# --- Remove the "24Slides" attribution picture (last shape on the slide) ---
$pic = $s.Shapes.Item($s.Shapes.Count)
$pic.Delete()
